# Weekly update: insert the new week's Apio (Americana) price rows at the
# top of the "Terminal La Palmera de La Serena" block (row 291), pushing the
# existing historical rows down by two rows (one row per quality grade:
# "Primera" and "Segunda"). The dimension grows from R312 to R314.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current first data row of the block.
$ws.Rows.Item(291).Insert()
$ws.Rows.Item(291).Insert()

# Columns that are constant across the whole "Apio" block.
$ws.Range("A291:A292").Value = 8
$ws.Range("B291:B292").Value = "Terminal La Palmera de La Serena"
$ws.Range("C291:C292").Value = "Coquimbo"
$ws.Range("E291:E292").Value = 4
$ws.Range("F291:F292").Value = 100112017
$ws.Range("G291:G292").Value = "Apio"
$ws.Range("H291:H292").Value = "Americana (o)"
$ws.Range("N291:N292").Value = "`$/docena de matas"
$ws.Range("O291:O292").Value = "Provincia del Elquí"
$ws.Range("Q291:Q292").Value = 6
$ws.Range("R291:R292").Value = "Hortaliza"

# Row 291: new week, quality "Primera".
$ws.Cells.Item(291, 4).Value = 44585
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 2400
$ws.Cells.Item(291, 11).Value = 8000
$ws.Cells.Item(291, 12).Value = 9000
$ws.Cells.Item(291, 13).Value = 8500
$ws.Cells.Item(291, 16).Value = 1417

# Row 292: new week, quality "Segunda".
$ws.Cells.Item(292, 4).Value = 44585
$ws.Cells.Item(292, 9).Value = "Segunda"
$ws.Cells.Item(292, 10).Value = 1400
$ws.Cells.Item(292, 11).Value = 6000
$ws.Cells.Item(292, 12).Value = 7000
$ws.Cells.Item(292, 13).Value = 6500
$ws.Cells.Item(292, 16).Value = 1083
